# QUERY_100set_1times.xlsx - "added clear cache before graph building and
# query performance test"
#
# The sheet used to hold:
#   A1 = "QUERY"                                          (shared string)
#   A2 = =QUERY("ShuffledDataSet100","ShuffledDataSet100_1_local")  (#NAME?)
#
# The commit removes the label in A1 and moves the QUERY(...) performance
# probe formula up into A1, leaving the sheet with a single formula cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "QUERY" label that lived in A1 (this is what empties
# sharedStrings.xml once nothing else references it).
$ws.Range("A1").ClearContents()

# Move the QUERY(...) performance-test formula from A2 up to A1.
$ws.Range("A1").Formula = '=QUERY("ShuffledDataSet100","ShuffledDataSet100_1_local")'

# A2 no longer holds anything.
$ws.Range("A2").ClearContents()

# Selection moves to the whole of row 1 (was a single cell, A3, previously).
$ws.Rows.Item(1).Select()
